$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.332.43"
$c.ClearFormats()
$ws.Range("E2").Value = "  -3.38%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.933.81"
$c.ClearFormats()
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("E4").Value = "  +0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "249.20"
$c.ClearFormats()
$ws.Range("E5").Value = "  -4.07%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.7230"
$c.ClearFormats()
$ws.Range("E6").Value = "  -8.87%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.23%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3288"
$c.ClearFormats()
$ws.Range("E8").Value = "  -8.90%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "27.62"
$c.ClearFormats()
$ws.Range("E9").Value = "  -3.64%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.06848"
$c.ClearFormats()
$ws.Range("E10").Value = "  -3.02%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.8072"
$c.ClearFormats()
$ws.Range("E11").Value = "  -4.69%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08066"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.33%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.933.36"
$c.ClearFormats()
$ws.Range("E13").Value = "  -3.59%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.415"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.93%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "94.91"
$c.ClearFormats()
$ws.Range("E15").Value = "  -6.47%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.51"
$c.ClearFormats()
$ws.Range("E16").Value = "  -1.03%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.328.61"
$c.ClearFormats()
$ws.Range("E17").Value = "  -3.29%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "251.07"
$c.ClearFormats()
$ws.Range("E18").Value = "  -9.08%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000008025"
$c.ClearFormats()
$ws.Range("E19").Value = "  +1.17%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "5.827"
$c.ClearFormats()
$ws.Range("E20").Value = "  -1.90%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.186.73"
$c.ClearFormats()
$ws.Range("E21").Value = "  -3.30%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.23%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.11%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.872"
$c.ClearFormats()
$ws.Range("E24").Value = "  -4.57%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.727"
$c.ClearFormats()
$ws.Range("E25").Value = "  -4.68%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "159.69"
$c.ClearFormats()
$ws.Range("E26").Value = "  -2.94%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.389"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  -4.65%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.1342"
$c.ClearFormats()
$ws.Range("E29").Value = "  -11.58%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.558"
$c.ClearFormats()
$ws.Range("E30").Value = "  -4.67%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.338"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.53%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.399"
$c.ClearFormats()
$ws.Range("E32").Value = "  -5.22%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.176"
$c.ClearFormats()
$ws.Range("E33").Value = "  -5.22%  "
$ws.Range("E34").Value = "  -2.53%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.220"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.39%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7405"
$c.ClearFormats()
$ws.Range("E36").Value = "  -3.49%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.752"
$c.ClearFormats()
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("E38").Value = "  -2.10%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.831"
$c.ClearFormats()
$ws.Range("E39").Value = "  -4.25%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.605"
$c.ClearFormats()
$ws.Range("E40").Value = "  -1.43%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "79.20"
$c.ClearFormats()
$ws.Range("E41").Value = "  -2.43%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.4465"
$c.ClearFormats()
$ws.Range("E42").Value = "  -5.87%  "
$ws.Range("E43").Value = "  -9.54%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.19%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.8353"
$c.ClearFormats()
$ws.Range("E45").Value = "  -2.79%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "102.07"
$c.ClearFormats()
$ws.Range("E46").Value = "  -2.69%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.292"
$c.ClearFormats()
$ws.Range("E48").Value = "  -4.93%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "36.42"
$c.ClearFormats()
$ws.Range("E49").Value = "  -1.52%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.05952"
$c.ClearFormats()
$ws.Range("E50").Value = "  -0.31%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.474"
$c.ClearFormats()
$ws.Range("E51").Value = "  -0.68%  "
